# Update LR-pair TPM values (Csf3-Csf3r) per new TPM-based script output.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value  = 0.03763333333333333   # G2
$ws.Cells.Item(2, 8).Value  = 0.1129                 # H2
$ws.Cells.Item(2, 9).Value  = 0.1820809155331632     # I2
$ws.Cells.Item(2, 10).Value = 0.1820809155331632     # J2
$ws.Cells.Item(2, 11).Value = 3                       # K2
$ws.Cells.Item(2, 12).Value = 1                       # L2
$ws.Cells.Item(2, 13).Value = 0.02144466666666667    # M2
$ws.Cells.Item(2, 14).Value = 0.064334               # N2
$ws.Cells.Item(2, 17).Value = 0.0008070342888888889  # Q2
$ws.Cells.Item(2, 18).Value = 0.007263308600000001   # R2
$ws.Cells.Item(2, 19).Value = 0.1820809155331632     # S2
$ws.Cells.Item(2, 20).Value = 0.1820809155331632     # T2

# Row 3
$ws.Cells.Item(3, 9).Value  = 0.5576820728517193     # I3
$ws.Cells.Item(3, 10).Value = 0.5576820728517194     # J3
$ws.Cells.Item(3, 11).Value = 3                       # K3
$ws.Cells.Item(3, 12).Value = 1                       # L3
$ws.Cells.Item(3, 13).Value = 0.02144466666666667    # M3
$ws.Cells.Item(3, 14).Value = 0.064334               # N3
$ws.Cells.Item(3, 17).Value = 0.002471805206888889   # Q3
$ws.Cells.Item(3, 18).Value = 0.022246246862         # R3
$ws.Cells.Item(3, 19).Value = 0.5576820728517193     # S3
$ws.Cells.Item(3, 20).Value = 0.5576820728517194     # T3

# Row 4
$ws.Cells.Item(4, 7).Value  = 0.053787               # G4
$ws.Cells.Item(4, 9).Value  = 0.2602370116151174     # I4
$ws.Cells.Item(4, 10).Value = 0.2602370116151174     # J4
$ws.Cells.Item(4, 11).Value = 3                       # K4
$ws.Cells.Item(4, 12).Value = 1                       # L4
$ws.Cells.Item(4, 13).Value = 0.02144466666666667    # M4
$ws.Cells.Item(4, 14).Value = 0.064334               # N4
$ws.Cells.Item(4, 17).Value = 0.001153444286         # Q4
$ws.Cells.Item(4, 18).Value = 0.010380998574         # R4
$ws.Cells.Item(4, 19).Value = 0.2602370116151174     # S4
$ws.Cells.Item(4, 20).Value = 0.2602370116151174     # T4
